$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New measurement data (resampled dataset, per commit "custom accuracy + 데이터 1000개")
$arr = New-Object 'object[,]' 4,34
$arr[0,0] = 45154.50694444445
$arr[0,1] = 10.726
$arr[0,2] = 7.333
$arr[0,3] = 3.404
$arr[0,4] = 23.56
$arr[0,5] = 17.15
$arr[0,6] = 8.176
$arr[0,7] = 24.228
$arr[0,8] = 13.347
$arr[0,9] = 5.245
$arr[0,10] = 7.323
$arr[0,11] = 9.308999999999999
$arr[0,12] = 10.191
$arr[0,13] = 2.44
$arr[0,14] = 8.647
$arr[0,15] = 11.655
$arr[0,16] = 7.955
$arr[0,17] = 2.648
$arr[0,18] = 1.093
$arr[0,19] = 124.223
$arr[0,20] = 23.834
$arr[0,21] = 7.982
$arr[0,22] = 14.964
$arr[0,23] = 8.048999999999999
$arr[0,24] = 2.19
$arr[0,25] = 13.597
$arr[0,26] = 7.05
$arr[0,27] = 6.629
$arr[0,28] = 7.562
$arr[0,29] = 9.942
$arr[0,30] = 2.682
$arr[0,31] = 21.628
$arr[0,32] = 4.075
$arr[0,33] = 9.978

$arr[1,0] = 45154.51388888889
$arr[1,1] = 23.804
$arr[1,2] = 17.612
$arr[1,3] = 2.007
$arr[1,4] = 52.167
$arr[1,5] = 41.902
$arr[1,6] = 18.61
$arr[1,7] = 69.377
$arr[1,8] = 29.055
$arr[1,9] = 12.799
$arr[1,10] = 18.592
$arr[1,11] = 20.894
$arr[1,12] = 22.327
$arr[1,13] = 5.817
$arr[1,14] = 18.798
$arr[1,15] = 26.547
$arr[1,16] = 16.108
$arr[1,17] = 1.4
$arr[1,18] = 0.993
$arr[1,19] = 278.812
$arr[1,20] = 52.521
$arr[1,21] = 17.351
$arr[1,22] = 34.973
$arr[1,23] = 18.397
$arr[1,24] = 2.874
$arr[1,25] = 35.068
$arr[1,26] = 15.326
$arr[1,27] = 13.719
$arr[1,28] = 16.078
$arr[1,29] = 21.989
$arr[1,30] = 1.136
$arr[1,31] = 63.278
$arr[1,32] = 9.619
$arr[1,33] = 21.693

$arr[2,0] = 45154.52083333334
$arr[2,1] = 1.278
$arr[2,2] = 0.8100000000000001
$arr[2,3] = 0.797
$arr[2,4] = 3.1
$arr[2,5] = 1.725
$arr[2,6] = 0.9360000000000001
$arr[2,7] = 13.257
$arr[2,8] = 1.716
$arr[2,9] = 0.827
$arr[2,10] = 0.524
$arr[2,11] = 1.227
$arr[2,12] = 1.463
$arr[2,13] = 0.237
$arr[2,14] = 1.128
$arr[2,15] = 1.628
$arr[2,16] = 1.321
$arr[2,17] = 0.8090000000000001
$arr[2,18] = 0.204
$arr[2,19] = 10.019
$arr[2,20] = 3.801
$arr[2,21] = 1.041
$arr[2,22] = 2.393
$arr[2,23] = 1.097
$arr[2,24] = 0.401
$arr[2,25] = 6.094
$arr[2,26] = 0.92
$arr[2,27] = 1.018
$arr[2,28] = 1.137
$arr[2,29] = 1.295
$arr[2,30] = 0.722
$arr[2,31] = 12.874
$arr[2,32] = 0.398
$arr[2,33] = 1.316

$arr[3,0] = 45154.52777777778
$arr[3,1] = 11.4
$arr[3,2] = 8.470000000000001
$arr[3,3] = 0.93
$arr[3,4] = 25.06
$arr[3,5] = 20.09
$arr[3,6] = 8.92
$arr[3,7] = 32.39
$arr[3,8] = 13.93
$arr[3,9] = 6.16
$arr[3,10] = 8.869999999999999
$arr[3,11] = 10.05
$arr[3,12] = 10.76
$arr[3,13] = 2.76
$arr[3,14] = 9.02
$arr[3,15] = 12.68
$arr[3,16] = 7.77
$arr[3,17] = 0.7
$arr[3,18] = 0.45
$arr[3,19] = 129.97
$arr[3,20] = 25.06
$arr[3,21] = 8.33
$arr[3,22] = 16.64
$arr[3,23] = 8.800000000000001
$arr[3,24] = 1.37
$arr[3,25] = 16.01
$arr[3,26] = 7.36
$arr[3,27] = 6.6
$arr[3,28] = 7.74
$arr[3,29] = 10.57
$arr[3,30] = 0.53
$arr[3,31] = 29.01
$arr[3,32] = 4.6
$arr[3,33] = 10.41

$ws.Range("A2:AH5").Value = $arr

# Remove the now-unused 6th data row (5 data rows -> 4 data rows)
$ws.Rows.Item(6).Delete()

# Column width tweaks for J, O, Q, V, X, AA, AB, AC, AH (7 -> 8 chars)
$ws.Range("J1").EntireColumn.ColumnWidth = 7.125
$ws.Range("O1").EntireColumn.ColumnWidth = 7.125
$ws.Range("Q1").EntireColumn.ColumnWidth = 7.125
$ws.Range("V1").EntireColumn.ColumnWidth = 7.125
$ws.Range("X1").EntireColumn.ColumnWidth = 7.125
$ws.Range("AA1").EntireColumn.ColumnWidth = 7.125
$ws.Range("AB1").EntireColumn.ColumnWidth = 7.125
$ws.Range("AC1").EntireColumn.ColumnWidth = 7.125
$ws.Range("AH1").EntireColumn.ColumnWidth = 7.125
